$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - set values then copy H1's format (style) onto them
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new columns I and J (rows 2-10)
$iValues = @(3, 9, 4, 7, 8, 5, 8, 5, 1)
$jValues = @(4, 9, 5, 7, 9, 6, 9, 7, 2)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
